$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) - Micro / SMEs
$ws.Range("B11").Value = "'63.96"
$ws.Range("C11").Value = "'5.34"

# Employment (% of total) - Micro / SMEs / MSMEs
$ws.Range("B12").Value = "'24.88"
$ws.Range("C12").Value = "'41.74"
$ws.Range("D12").Value = "'66.62"

# Enterprises (% of total) - Micro / SMEs / MSMEs
$ws.Range("B14").Value = "'87.05"
$ws.Range("C14").Value = "'12.58"
$ws.Range("D14").Value = "'99.63"

# Enterprises density (per 1000 people) - second table (SME Associations)
$ws.Range("B37").Value = "'11.78"
$ws.Range("C37").Value = "'5.12"
$ws.Range("D37").Value = "'16.89"

# Employment (% of total) - second table (SME Associations)
$ws.Range("B38").Value = "'16.14"
$ws.Range("C38").Value = "'50.51"
$ws.Range("D38").Value = "'66.65"

# Value added to the economy (% of total) - SBS Eurostat
$ws.Range("B44").Value = "'11.32"
$ws.Range("C44").Value = "'47.55"
$ws.Range("D44").Value = "'58.88"
